$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column C (y_1) and column E (y_1_forecast), rows 2-19
$values = @{
    2  = @{ C = 1.981300392979479;   E = 2.228397109637226 }
    3  = @{ C = -0.4626567965509865; E = -0.2018858887078534 }
    4  = @{ C = -0.4861108058567654; E = 0.2108047537406676 }
    5  = @{ C = 2.490556233265284;   E = 0.9684279156219944 }
    6  = @{ C = 1.175152734113438;   E = 1.69818237209749 }
    7  = @{ C = -0.2674335569108899; E = 0.679606031449409 }
    8  = @{ C = 2.038609866767938;   E = 1.325176859452326 }
    9  = @{ C = 1.566972126803345;   E = 1.522808462763714 }
    10 = @{ C = 2.246337373619012;   E = 1.693557061600992 }
    11 = @{ C = 1.777150434343522;   E = 1.905564797014625 }
    12 = @{ C = 1.741137453897301;   E = 2.082477074609068 }
    13 = @{ C = 1.562095320687407;   E = 1.845103901518885 }
    14 = @{ C = -0.8985012482809474; E = 0.02570757229449772 }
    15 = @{ C = -0.7241284594088016; E = -0.1006764062950749 }
    16 = @{ C = 3.72986100763808;    E = 1.357201584213352 }
    17 = @{ C = -0.3358560228777674; E = 1.155533515694374 }
    18 = @{ C = -1.054811008161194;  E = -0.009717596728553435 }
    19 = @{ C = 1.593155398714674;   E = 0.1924237296221154 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}
